$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.370.48"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -2.07%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.710.33"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.27"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.77%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5341"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.14%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2665"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06601"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.87"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07629"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.573"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.704.02"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.46%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.946.15"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.93%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5750"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8174"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.92%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.80"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.355.01"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.12%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.62"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.676"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -3.68%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.47"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -4.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.978"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.13%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "141.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -2.74%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.737"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1215"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.285"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.68%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.33"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05408"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.295"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.73%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.501"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -5.61%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.430"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.86%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.647"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.49%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.880"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.58%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9504"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.417"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.19%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5869"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01633"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.37%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.869"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.044.47"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8421"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.94"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.28%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.852.11"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.99%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.67%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.96%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4509"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.59%  "

$ws.Range("E49").Value = "  +0.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.108"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.51%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05240"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.86%  "
